# Auto-generated Excel COM-interop script to update the cryptos worksheet
# with the latest scraped prices / 1h volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D holds numeric-looking text (e.g. "29.212.32", "0.07746") that
# --- must stay plain text (it is not a real number - some values use dots as
# --- thousands separators, others would lose significant trailing zeros).
# --- Force text format before writing, then restore the default "Normal" style
# --- so no stray formatting is left behind on the cells.
$dRanges = @(
    "D2:D26",
    "D28:D29",
    "D31:D51",
)
foreach ($addr in $dRanges) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "29.212.32"
$ws.Range("E2").Value = "  -1.89%  "

# Row 3
$ws.Range("D3").Value = "1.858.38"
$ws.Range("E3").Value = "  -0.72%  "

# Row 4
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "0.6922"
$ws.Range("E5").Value = "  -3.49%  "

# Row 6
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "238.44"
$ws.Range("E6").Value = "  -1.32%  "

# Row 7
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  -0.25%  "

# Row 8
$ws.Range("D8").Value = "0.07746"
$ws.Range("E8").Value = "  +2.92%  "

# Row 9
$ws.Range("D9").Value = "0.3059"
$ws.Range("E9").Value = "  -2.77%  "

# Row 10
$ws.Range("D10").Value = "23.35"
$ws.Range("E10").Value = "  -4.71%  "

# Row 11
$ws.Range("D11").Value = "0.08061"
$ws.Range("E11").Value = "  -1.53%  "

# Row 12
$ws.Range("D12").Value = "1.897.07"
$ws.Range("E12").Value = "  +0.96%  "

# Row 13
$ws.Range("D13").Value = "0.7247"
$ws.Range("E13").Value = "  -2.40%  "

# Row 14
$ws.Range("D14").Value = "5.218"
$ws.Range("E14").Value = "  -1.99%  "

# Row 15
$ws.Range("D15").Value = "89.59"
$ws.Range("E15").Value = "  -3.03%  "

# Row 16
$ws.Range("D16").Value = "29.212.55"
$ws.Range("E16").Value = "  -1.72%  "

# Row 17
$ws.Range("D17").Value = "5.757"
$ws.Range("E17").Value = "  -4.12%  "

# Row 18
$ws.Range("D18").Value = "0.000007823"
$ws.Range("E18").Value = "  -1.18%  "

# Row 19
$ws.Range("D19").Value = "13.27"
$ws.Range("E19").Value = "  -1.30%  "

# Row 20
$ws.Range("D20").Value = "235.75"
$ws.Range("E20").Value = "  -4.24%  "

# Row 21
$ws.Range("D21").Value = "0.9986"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22
$ws.Range("D22").Value = "2.104.51"
$ws.Range("E22").Value = "  -0.20%  "

# Row 23
$ws.Range("D23").Value = "0.9987"
$ws.Range("E23").Value = "  -0.36%  "

# Row 24
$ws.Range("D24").Value = "7.475"
$ws.Range("E24").Value = "  -2.95%  "

# Row 25
$ws.Range("D25").Value = "162.10"
$ws.Range("E25").Value = "  -0.91%  "

# Row 26
$ws.Range("D26").Value = "8.989"
$ws.Range("E26").Value = "  -2.01%  "

# Row 27
$ws.Range("E27").Value = "  -3.19%  "

# Row 28
$ws.Range("D28").Value = "18.11"
$ws.Range("E28").Value = "  -2.30%  "

# Row 29
$ws.Range("D29").Value = "1.965"
$ws.Range("E29").Value = "  -1.79%  "

# Row 30
$ws.Range("E30").Value = "  -1.39%  "

# Row 31
$ws.Range("D31").Value = "4.522"
$ws.Range("E31").Value = "  -0.27%  "

# Row 32
$ws.Range("D32").Value = "1.488"
$ws.Range("E32").Value = "  -2.38%  "

# Row 33
$ws.Range("D33").Value = "4.028"
$ws.Range("E33").Value = "  -3.49%  "

# Row 34
$ws.Range("D34").Value = "0.05195"
$ws.Range("E34").Value = "  -4.65%  "

# Row 35
$ws.Range("D35").Value = "1.188"
$ws.Range("E35").Value = "  -2.77%  "

# Row 36
$ws.Range("D36").Value = "0.7065"
$ws.Range("E36").Value = "  -4.02%  "

# Row 37
$ws.Range("D37").Value = "1.020"
$ws.Range("E37").Value = "  +2.14%  "

# Row 38
$ws.Range("D38").Value = "2.669"
$ws.Range("E38").Value = "  -1.25%  "

# Row 39
$ws.Range("D39").Value = "0.01856"
$ws.Range("E39").Value = "  -2.71%  "

# Row 40
$ws.Range("D40").Value = "2.678"
$ws.Range("E40").Value = "  -2.04%  "

# Row 41
$ws.Range("D41").Value = "0.9247"
$ws.Range("E41").Value = "  +4.84%  "

# Row 42
$ws.Range("D42").Value = "1.100.90"
$ws.Range("E42").Value = "  +6.34%  "

# Row 43
$ws.Range("D43").Value = "5.973"
$ws.Range("E43").Value = "  -0.17%  "

# Row 44
$ws.Range("D44").Value = "0.4296"
$ws.Range("E44").Value = "  -3.39%  "

# Row 45
$ws.Range("D45").Value = "70.86"
$ws.Range("E45").Value = "  -0.67%  "

# Row 46
$ws.Range("D46").Value = "0.9988"
$ws.Range("E46").Value = "  -0.33%  "

# Row 47
$ws.Range("D47").Value = "102.22"
$ws.Range("E47").Value = "  -1.31%  "

# Row 48
$ws.Range("D48").Value = "1.795"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49
$ws.Range("D49").Value = "2.001.67"
$ws.Range("E49").Value = "  -0.01%  "

# Row 50
$ws.Range("D50").Value = "9.200"
$ws.Range("E50").Value = "  -3.78%  "

# Row 51
$ws.Range("D51").Value = "7.027"
$ws.Range("E51").Value = "  -5.72%  "

# Restore default styling on the D column cells we reformatted above so the
# saved workbook has no leftover style deltas versus the original.
foreach ($addr in $dRanges) {
    $ws.Range($addr).Style = "Normal"
}

Write-Output "cryptos worksheet updated"
